$d = $word.ActiveDocument

# 1. "K-8 Intro to Computer Science" -> "Computer Science Fundamentals"
#    (affects both the bold title-line mention and the "WHAT IS IT?" body paragraph)
$d.Content.Find.Execute("K-8 Intro to Computer Science", $true, $false, $false, $false, $false, $true, 1, $false, "Computer Science Fundamentals", 2)

# 2. "90% of American schools don" -> "Most American schools don"
$d.Content.Find.Execute("90% of American schools don", $true, $false, $false, $false, $false, $true, 1, $false, "Most American schools don", 2)

# 3. "Fewer students are learning how computers work than a decade ago!" -> "Computing jobs are the number 1 source of new wages in the US!"
$d.Content.Find.Execute("Fewer students are learning how computers work than a decade ago!", $true, $false, $false, $false, $false, $true, 1, $false, "Computing jobs are the number 1 source of new wages in the US!", 2)

# 4. "irls and students of color are severely underrepresented in computer science." -> "irls and many minorities are severely underrepresented in computer science."
$d.Content.Find.Execute("irls and students of color are severely underrepresented in computer science.", $true, $false, $false, $false, $false, $true, 1, $false, "irls and many minorities are severely underrepresented in computer science.", 2)

# 5. "...girls and students of color. " -> "...girls and underrepresented minorities. " (Code.org blurb near the end)
$d.Content.Find.Execute("expanding participation by girls and students of color. ", $true, $false, $false, $false, $false, $true, 1, $false, "expanding participation by girls and underrepresented minorities. ", 2)

# 6. Move the _GoBack bookmark from the paragraph right after the banner image
#    down to the empty paragraph right after "...severely underrepresented in computer science."
foreach ($bm in @($d.Bookmarks)) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "severely underrepresented in computer science\.") {
        $target = $p
    }
}
$nextPara = $target.Next()
$d.Bookmarks.Add("_GoBack", $nextPara.Range)
